$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 241, shifting existing rows 241:261 down to 242:262
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new data record
$ws.Cells.Item(241, 1).Value = 3
$ws.Cells.Item(241, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(241, 3).Value = "Coquimbo"
$ws.Cells.Item(241, 4).Value = 45223
$ws.Cells.Item(241, 5).Value = 5
$ws.Cells.Item(241, 6).Value = 100112052
$ws.Cells.Item(241, 7).Value = "Albahaca"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 60
$ws.Cells.Item(241, 11).Value = 5000
$ws.Cells.Item(241, 12).Value = 5000
$ws.Cells.Item(241, 13).Value = 5000
$ws.Cells.Item(241, 14).Value = "$/docena de matas"
$ws.Cells.Item(241, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(241, 16).Value = 833
$ws.Cells.Item(241, 17).Value = 6
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# Ensure the date cell keeps the workbook's date number format (matches column D elsewhere)
$ws.Cells.Item(241, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
